# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> used by the slide master (was "Integral")
#   ppt/theme/theme2.xml -> used by the notes master (was "Office Theme")
# The authored change swaps the two themes' contents: the slide master
# ends up with the "Office Theme" palette, and the notes master ends up
# with the "Integral" palette. The font scheme and format scheme are
# identical between the two themes (both "Office"/Arial), so the only
# substantive difference is the 12 theme colours.
#
# Apply it through the Office colour-scheme object model, which edits
# the live theme colours (and is reflected back into the underlying
# theme XML) rather than poking at file internals directly.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Target palette = the "Office Theme" colours, in
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order.
# RGB(r,g,b) -> r + g*256 + b*65536 (the OLE_COLOR encoding PowerPoint uses).
$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72
